# Myst.pptx edit script
# 1) Slide 2 ("Il problema") body text box: the two body paragraphs each had
#    their text split across multiple <a:r> runs (because of incremental
#    edits in the authoring tool). The runs are collapsed back into a single
#    run per paragraph, keeping the very first run's formatting, by
#    re-assigning the full paragraph text through a Characters() range (this
#    forces PowerPoint to rebuild the paragraph as one run instead of
#    patching the existing runs in place).
# 2) Slide 3 ("La soluzione") body text box: first paragraph's wording
#    changes from "un sito web" to "una PWA", and the second paragraph's
#    seven runs are likewise collapsed into a single run.

function Set-ParagraphText {
    param($Shape, $ParaIndex, $NewText)
    $tr = $Shape.TextFrame.TextRange
    $para = $tr.Paragraphs($ParaIndex, 1)
    $start = $para.Start
    $len = $para.Length
    $oldText = $para.Text
    # A non-final paragraph's range includes a trailing paragraph-mark
    # character (\r); the very last paragraph in the text frame does not.
    # Only trim the range length when that trailing \r is actually present,
    # so we replace just the visible text and leave the paragraph mark as-is.
    if ($len -gt 0 -and [int][char]$oldText.Substring($oldText.Length - 1) -eq 13) {
        $len = $len - 1
    }
    $chars = $tr.Characters($start, $len)
    $chars.Text = $NewText
}

$p = $ppt.ActivePresentation

# --- Slide 2: "Il problema" --------------------------------------------
$slide2 = $p.Slides.Item(2)
$body2 = $slide2.Shapes.Item(2)

Set-ParagraphText $body2 1 "Rispetto al mercato fisico, il mercato digitale dei videogiochi non permette la rivendita o lo scambio della licenza dei propri giochi, non permettendo alle persone di risparmiare."

Set-ParagraphText $body2 2 "Anche se è possibile pagare di meno tramite le chiavi di attivazione di giochi esse rappresentano una zona grigia del mercato digitale a causa della loro dubbia provenienza in alcuni casi."

# --- Slide 3: "La soluzione" --------------------------------------------
$slide3 = $p.Slides.Item(3)
$body3 = $slide3.Shapes.Item(2)

Set-ParagraphText $body3 1 "Si vuole sviluppare una PWA per una piattaforma che vende videogiochi digitali."

Set-ParagraphText $body3 2 "L’utente può registrare un account, comprare giochi  (che verranno aggiunti alla sua libreria) e ha la possibilità di rivendere poi i giochi comprati ad altri utenti(aggiustati per valuta in base al paese dell’utente)."
